$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New aggregated counts query (Programs/Studies/Cases/Samples/Case Files/Study Files)
# used by CasesTab, SamplesTab and FilesTab dbExcel (column C) queries
$newAggQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE samp.summarized_sample_type IN ["Normal Tissue"]  
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Updated File query (FilesTab WebExcel / column B): drop the trailing Study Code line
$newFileQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
WHERE samp.summarized_sample_type IN ["Normal Tissue"] 
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis  
'@

# Updated Case query (CasesTab WebExcel / column B): append the Cohort column
$newCaseQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
WHERE samp.summarized_sample_type IN ["Normal Tissue"] 
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

# Populate the dbExcel (column C) queries first so the new shared string is
# created in the same relative order the workbook ends up with.
$ws.Range("C2").Value = $newAggQuery
$ws.Range("C3").Value = $newAggQuery
$ws.Range("C4").Value = $newAggQuery

# FilesTab WebExcel query (column B)
$ws.Range("B4").Value = $newFileQuery

# CasesTab WebExcel query (column B) - Samples (B3) stays as-is
$ws.Range("B2").Value = $newCaseQuery

# Row heights shrink now that the stored text is shorter
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# Selection moves from B3 to F3 and the view scrolls back to the top-left
$ws.Range("F3").Select() | Out-Null
